$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S2").Value = 6
$ws.Range("S3").Value = 9.75
$ws.Range("S7").Value = 8.6999999999999993
$ws.Range("S9").Value = 8.4
$ws.Range("R10").Value = 6.75
$ws.Range("S10").Value = 4
$ws.Range("R11").Value = 9
$ws.Range("S11").Value = 5.8
$ws.Range("S12").Value = 5.0999999999999996
$ws.Range("S14").Value = 5.2
$ws.Range("S15").Value = 9
$ws.Range("R16").Value = 10
$ws.Range("S17").Value = 7.1
$ws.Range("R20").Value = 9
$ws.Range("S20").Value = 5.25
$ws.Range("S21").Value = 8.3000000000000007
$ws.Range("R22").Value = 9
$ws.Range("S22").Value = 4.5999999999999996

$ws.Range("Q2:U23").Select()
